$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q3").Value = 1.83
$ws.Range("H4").Value = 3
$ws.Range("H5").Value = 2.28
$ws.Range("I5").Value = 2.72
$ws.Range("Q5").Value = 1.74
$ws.Range("G7").Value = 1.49
$ws.Range("G11").Value = 2.44
$ws.Range("F12").Value = 1.51
$ws.Range("G12").Value = 1.55
$ws.Range("H12").Value = 8.4
$ws.Range("I12").Value = 9.8
